$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.053392887115479
$ws.Range("B1").Value = 4.300155162811279
$ws.Range("C1").Value = 1.065963387489319
$ws.Range("D1").Value = 0.2835461795330048
$ws.Range("E1").Value = 0.1828610301017761
